$wb = $excel.ActiveWorkbook

# --- Keywords sheet: zoom out from 175% to 102% (it was the active tab, but won't be anymore) ---
$wsKeywords = $wb.Worksheets.Item("Keywords")
$wsKeywords.Activate()
$wsKeywords.Range("A3").Select()
$excel.ActiveWindow.Zoom = 102

# --- ColumnHeaders sheet: update the missing-value explanation for station_distance (G18) ---
$wsHeaders = $wb.Worksheets.Item("ColumnHeaders")
$wsHeaders.Range("G18").Value = "Distance to nearest station is greater than 2 km"

# Make ColumnHeaders the active sheet/tab with G18 selected
$wsHeaders.Activate()
$wsHeaders.Range("G18").Select()
